# Update the yearly database: roll the 12-month trailing window forward one
# year (drop 1396/12, shift 1397/12..1400/12 left, add 1401/12) and refresh
# the underlying figures per the "change read_price algorithm" update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 header labels (period headers) - shift one year forward
$ws.Range("E8").Value = 'دوازده ماهه منتهی به 1397/12'
$ws.Range("F8").Value = 'دوازده ماهه منتهی به 1398/12'
$ws.Range("G8").Value = 'دوازده ماهه منتهی به 1399/12'
$ws.Range("H8").Value = 'دوازده ماهه منتهی به 1400/12'
$ws.Range("I8").Value = 'دوازده ماهه منتهی به 1401/12'

# General & administrative expense rows - shift figures one year forward and
# append the new 1401/12 column value
$ws.Range("G10").Value = 52391
$ws.Range("H10").Value = 18178
$ws.Range("I10").Value = 18918
$ws.Range("E12").Value = 32663
$ws.Range("F12").Value = 50134
$ws.Range("G12").Value = 24113
$ws.Range("H12").Value = 0
$ws.Range("E13").Value = 1658
$ws.Range("F13").Value = 1298
$ws.Range("G13").Value = 1434
$ws.Range("H13").Value = 1029
$ws.Range("I13").Value = 3245
$ws.Range("E15").Value = 979
$ws.Range("F15").Value = 967
$ws.Range("G15").Value = 856
$ws.Range("H15").Value = 3062
$ws.Range("I15").Value = 3251
$ws.Range("E16").Value = 984
$ws.Range("F16").Value = 1962
$ws.Range("G16").Value = 2298
$ws.Range("H16").Value = 2344
$ws.Range("I16").Value = 4145
$ws.Range("E17").Value = 68294
$ws.Range("F17").Value = 110813
$ws.Range("G17").Value = 148320
$ws.Range("H17").Value = 258583
$ws.Range("I17").Value = 407601
$ws.Range("E19").Value = 63893
$ws.Range("F19").Value = 136873
$ws.Range("G19").Value = 117673
$ws.Range("H19").Value = 155871
$ws.Range("I19").Value = 306862
$ws.Range("E20").Value = 168471
$ws.Range("F20").Value = 302047
$ws.Range("G20").Value = 347085
$ws.Range("H20").Value = 439067
$ws.Range("I20").Value = 744022

# Row 24 header labels (personnel count table) - shift one year forward
$ws.Range("E24").Value = 'دوازده ماهه منتهی به 1397/12'
$ws.Range("F24").Value = 'دوازده ماهه منتهی به 1398/12'
$ws.Range("G24").Value = 'دوازده ماهه منتهی به 1399/12'
$ws.Range("H24").Value = 'دوازده ماهه منتهی به 1400/12'
$ws.Range("I24").Value = 'دوازده ماهه منتهی به 1401/12'

# Personnel count rows - shift figures one year forward and append new value
$ws.Range("E26").Value = 91
$ws.Range("F26").Value = 87
$ws.Range("G26").Value = 85
$ws.Range("I26").Value = 84
$ws.Range("E27").Value = 795
$ws.Range("F27").Value = 775
$ws.Range("G27").Value = 755
$ws.Range("H27").Value = 701
$ws.Range("I27").Value = 667
